# "Search channel from discover changes"
# - Login sheet gets new credentials in row 2 (manisha.kisan17@gmail.com / manisha123)
# - The old credentials (exh_0006@mailinator.com / chan@kisan18) move down into row 3,
#   styled with the blue bold "Calibri 13" font and re-hyperlinked.
# - Login becomes the active sheet/tab; Channel's selection moves to A4 and loses focus.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsChannel = $wb.Worksheets.Item("Channel")

# --- Row 2: swap in the new exhibitor credentials -------------------------
$wsLogin.Range("E2").Value = "manisha.kisan17@gmail.com"
$wsLogin.Range("F2").Value = "manisha123"

# --- Row 3: carry the former credentials down, with the new blue/bold style
$wsLogin.Range("E3").Value = "exh_0006@mailinator.com"
$wsLogin.Range("E3").Font.Color = 16711680
$wsLogin.Range("F3").Value = "chan@kisan18"
$wsLogin.Range("F3").Font.Color = 16711680

# --- Rebuild the hyperlinks (Login sheet) so each cell's link/display text
# matches its current value; clearing first avoids stacking stale entries.
$wsLogin.Hyperlinks.Delete() | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("E2"), "mailto:manisha.kisan17@gmail.com", [Type]::Missing, [Type]::Missing, "manisha.kisan17@gmail.com") | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("E3"), "mailto:exh_0006@mailinator.com", [Type]::Missing, [Type]::Missing, "exh_0006@mailinator.com") | Out-Null
$wsLogin.Hyperlinks.Add($wsLogin.Range("F3"), "chan@kisan18", [Type]::Missing, [Type]::Missing, "chan@kisan18") | Out-Null

# --- View state: Channel's cursor moves to A4 and it loses the active tab,
# Login's cursor moves to F14 and it becomes the active tab (selected last).
$wsChannel.Range("A4").Select() | Out-Null
$wsLogin.Range("F14").Select() | Out-Null
